# QA Compiler: Fix Actual Issues % negative value bug
#
# The "Actual Issues" percentage columns were showing "100.0%" but should
# show "100%" (the underlying fix clamps the computed percentage to the
# 0-100% range; these cached display strings just need to match the new
# formatted output). This updates the 11 affected cells across the DAILY
# and TOTAL sheets.
#
# Because these cells hold plain text that merely *looks* like a number,
# a bare assignment such as  $range.Value = "100%"  would make Excel's
# input parser treat it as the numeric value 1 formatted as a percentage
# (changing the cell's type/format). Prefixing the literal with a single
# quote tells Excel to store it verbatim as text, exactly like the
# original "100.0%" entries were stored, while leaving every other
# attribute of the cell (its style/border/fill) untouched.

$wb = $excel.ActiveWorkbook

$daily = $wb.Worksheets.Item("DAILY")
$dailyCells = @("E5", "I6", "M6", "E8", "I8", "M8")
foreach ($ref in $dailyCells) {
    $daily.Range($ref).Value = "'100%"
}

$total = $wb.Worksheets.Item("TOTAL")
$totalCells = @("C3", "C4", "C5", "C6", "C9")
foreach ($ref in $totalCells) {
    $total.Range($ref).Value = "'100%"
}
